{"js": "// Replace each two-digit multiplication problem's text with its new value.\n// Each original value is unique within the document, so a plain text search\n// safely targets exactly one run.\nconst replacements = [\n  [\"34\u00d730=\", \"66\u00d754=\"],\n  [\"99\u00d730=\", \"70\u00d752=\"],\n  [\"77\u00d798=\", \"48\u00d768=\"],\n  [\"46\u00d731=\", \"60\u00d757=\"],\n  [\"23\u00d787=\", \"88\u00d794=\"],\n  [\"36\u00d798=\", \"67\u00d771=\"],\n  [\"11\u00d720=\", \"75\u00d740=\"],\n  [\"67\u00d751=\", \"67\u00d789=\"],\n  [\"76\u00d737=\", \"31\u00d764=\"],\n  [\"36\u00d759=\", \"53\u00d731=\"],\n  [\"51\u00d778=\", \"80\u00d762=\"],\n  [\"48\u00d770=\", \"37\u00d732=\"],\n  [\"83\u00d749=\", \"50\u00d718=\"],\n  [\"96\u00d716=\", \"13\u00d730=\"],\n  [\"23\u00d711=\", \"30\u00d717=\"],\n  [\"56\u00d730=\", \"54\u00d719=\"],\n  [\"22\u00d786=\", \"74\u00d724=\"],\n  [\"47\u00d730=\", \"44\u00d769=\"],\n  [\"58\u00d784=\", \"56\u00d739=\"],\n  [\"58\u00d799=\", \"24\u00d724=\"],\n  [\"87\u00d761=\", \"19\u00d724=\"],\n  [\"77\u00d773=\", \"20\u00d799=\"],\n  [\"61\u00d755=\", \"20\u00d774=\"],\n  [\"76\u00d747=\", \"32\u00d722=\"],\n  [\"28\u00d796=\", \"30\u00d771=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication problem's text with its new value.\n# Each original value is unique within the document, so Find/Replace safely\n# targets exactly one run per pair.\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$replacements = @(\n  @(\"34\u00d730=\", \"66\u00d754=\"),\n  @(\"99\u00d730=\", \"70\u00d752=\"),\n  @(\"77\u00d798=\", \"48\u00d768=\"),\n  @(\"46\u00d731=\", \"60\u00d757=\"),\n  @(\"23\u00d787=\", \"88\u00d794=\"),\n  @(\"36\u00d798=\", \"67\u00d771=\"),\n  @(\"11\u00d720=\", \"75\u00d740=\"),\n  @(\"67\u00d751=\", \"67\u00d789=\"),\n  @(\"76\u00d737=\", \"31\u00d764=\"),\n  @(\"36\u00d759=\", \"53\u00d731=\"),\n  @(\"51\u00d778=\", \"80\u00d762=\"),\n  @(\"48\u00d770=\", \"37\u00d732=\"),\n  @(\"83\u00d749=\", \"50\u00d718=\"),\n  @(\"96\u00d716=\", \"13\u00d730=\"),\n  @(\"23\u00d711=\", \"30\u00d717=\"),\n  @(\"56\u00d730=\", \"54\u00d719=\"),\n  @(\"22\u00d786=\", \"74\u00d724=\"),\n  @(\"47\u00d730=\", \"44\u00d769=\"),\n  @(\"58\u00d784=\", \"56\u00d739=\"),\n  @(\"58\u00d799=\", \"24\u00d724=\"),\n  @(\"87\u00d761=\", \"19\u00d724=\"),\n  @(\"77\u00d773=\", \"20\u00d799=\"),\n  @(\"61\u00d755=\", \"20\u00d774=\"),\n  @(\"76\u00d747=\", \"32\u00d722=\"),\n  @(\"28\u00d796=\", \"30\u00d771=\"),\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n}\n"}
